$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")
$ws.Activate()

# Insert a new column before column N (14), shifting existing N..AO content right by one
$ws.Columns("N:N").Insert()

# New column header (row 1) -- "GiftCard2"
$ws.Cells.Item(1, 14).Value = "GiftCard2"

# New gift card code on the "Giftcard" test-data row (row 34); leading apostrophe
# mirrors how this alphanumeric code was originally typed so it keeps the
# text/quote-prefix formatting used by the neighbouring cell.
$ws.Cells.Item(34, 14).Value = "'2MT27C26L7277W77H44E"

# Leave the view roughly where the author left it after adding the column.
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("M16").Select()
